$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price column keeps its original text formatting so Excel
# does not silently re-interpret strings like "1.008" as numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "20.498.04"
$ws.Range("E2").Value = "  +1.98%  "
$ws.Range("D3").Value = "1.474.74"
$ws.Range("E3").Value = "  +3.45%  "
$ws.Range("D4").Value = "1.008"
$ws.Range("E4").Value = "  +1.09%  "
$ws.Range("D5").Value = "0.9610"
$ws.Range("E5").Value = "  -3.53%  "
$ws.Range("D6").Value = "276.49"
$ws.Range("E6").Value = "  -0.08%  "
$ws.Range("D7").Value = "0.3644"
$ws.Range("E7").Value = "  -1.52%  "
$ws.Range("D8").Value = "0.3042"
$ws.Range("E8").Value = "  -2.83%  "
$ws.Range("D9").Value = "39.65"
$ws.Range("E9").Value = "  -0.31%  "
$ws.Range("D10").Value = "1.048"
$ws.Range("E10").Value = "  -0.74%  "
$ws.Range("D11").Value = "0.06595"
$ws.Range("E11").Value = "  +0.71%  "
$ws.Range("E12").Value = "  +0.63%  "
$ws.Range("D13").Value = "18.11"
$ws.Range("E13").Value = "  +1.62%  "
$ws.Range("D14").Value = "5.449"
$ws.Range("E14").Value = "  -1.24%  "
$ws.Range("D15").Value = "6.158"
$ws.Range("E15").Value = "  -0.85%  "
$ws.Range("D16").Value = "0.00001027"
$ws.Range("E16").Value = "  +0.28%  "
$ws.Range("D17").Value = "1.478.45"
$ws.Range("E17").Value = "  +3.93%  "
$ws.Range("D18").Value = "0.05888"
$ws.Range("E18").Value = "  +3.55%  "
$ws.Range("D19").Value = "0.9675"
$ws.Range("E19").Value = "  -2.90%  "
$ws.Range("D20").Value = "69.04"
$ws.Range("E20").Value = "  -3.47%  "
$ws.Range("D21").Value = "5.458"
$ws.Range("E21").Value = "  -2.86%  "
$ws.Range("D22").Value = "14.45"
$ws.Range("E22").Value = "  -2.80%  "
$ws.Range("E23").Value = "  -1.09%  "
$ws.Range("D24").Value = "2.245"
$ws.Range("E24").Value = "  -0.24%  "
$ws.Range("D25").Value = "20.556.49"
$ws.Range("E25").Value = "  +2.28%  "
$ws.Range("D26").Value = "141.21"
$ws.Range("E26").Value = "  +5.54%  "
$ws.Range("D27").Value = "2.121"
$ws.Range("E27").Value = "  -7.58%  "
$ws.Range("D28").Value = "17.16"
$ws.Range("E28").Value = "  -0.82%  "
$ws.Range("D29").Value = "1.634.66"
$ws.Range("E29").Value = "  +3.58%  "
$ws.Range("D30").Value = "113.26"
$ws.Range("E30").Value = "  +2.33%  "
$ws.Range("D31").Value = "3.878"
$ws.Range("E31").Value = "  -0.92%  "
$ws.Range("D32").Value = "4.955"
$ws.Range("E32").Value = "  -6.18%  "
$ws.Range("D33").Value = "0.8085"
$ws.Range("E33").Value = "  -2.07%  "
$ws.Range("D34").Value = "0.07888"
$ws.Range("E34").Value = "  +1.11%  "
$ws.Range("D35").Value = "1.507"
$ws.Range("E35").Value = "  +2.85%  "
$ws.Range("D36").Value = "1.238"
$ws.Range("E36").Value = "  +11.69%  "
$ws.Range("D37").Value = "0.05752"
$ws.Range("E37").Value = "  -1.94%  "
$ws.Range("D38").Value = "4.721"
$ws.Range("E38").Value = "  -4.35%  "
$ws.Range("D39").Value = "7.694"
$ws.Range("E39").Value = "  -4.55%  "
$ws.Range("D40").Value = "0.02034"
$ws.Range("E40").Value = "  -1.46%  "
$ws.Range("D41").Value = "0.9614"
$ws.Range("E41").Value = "  -3.64%  "
$ws.Range("D42").Value = "10.40"
$ws.Range("E42").Value = "  -2.29%  "
$ws.Range("E43").Value = "  -0.40%  "
$ws.Range("D44").Value = "0.5265"
$ws.Range("E44").Value = "  -1.44%  "
$ws.Range("D45").Value = "3.503"
$ws.Range("E45").Value = "  -1.08%  "
$ws.Range("D46").Value = "12.03"
$ws.Range("E46").Value = "  -2.81%  "
$ws.Range("D47").Value = "116.87"
$ws.Range("E47").Value = "  -0.73%  "
$ws.Range("D48").Value = "0.5157"
$ws.Range("E48").Value = "  -1.38%  "
$ws.Range("E49").Value = "  -0.64%  "
$ws.Range("D50").Value = "0.06449"
$ws.Range("E50").Value = "  +3.48%  "
$ws.Range("D51").Value = "0.9929"
$ws.Range("E51").Value = "  -0.39%  "
